$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1").Value = "部位"
$ws.Range("J2").Value = "胸部"

$ws.Range("J2").Select()
